$wb = $excel.ActiveWorkbook
$wsTables = $wb.Worksheets.Item("Tables")
$wsFields = $wb.Worksheets.Item("Fields")

# ---------------------------------------------------------------------------
# Fields sheet: D18 (output_id / "Possible values") -> new plain-text string
# ---------------------------------------------------------------------------
$wsFields.Range("D18").Value = "any integer that corresponds to a unique output_id specified in the outputs table"

# ---------------------------------------------------------------------------
# Fields sheet: D19 (model_output_values / output_name / "Possible values")
# -> new rich text: bold "One of:" + newline + italic list of outputs
# ---------------------------------------------------------------------------
$fullText = "One of:`ncumulative fatalities`ncumulative cases`nICU beds per day`nventilators per day`nhospital admissions per day`nICU admissions per day`nfatalities per day`n(additional outputs to be added, in progress)"

$cell19 = $wsFields.Range("D19")
$cell19.Value = $fullText

# Touch the whole-cell Bold/Italic state (and revert) so that matching fonts
# get registered in the workbook's font table, mirroring how Excel tracks
# recently used fonts even though the cell itself keeps its original style.
$cell19.Font.Bold = $true
$cell19.Font.Bold = $false
$cell19.Font.Italic = $true
$cell19.Font.Italic = $false

$cell19.Characters(1,7).Font.Bold = $true
$cell19.Characters(9,191).Font.Italic = $true

# Row 19 grows taller to fit the new text
$wsFields.Rows.Item(19).RowHeight = 146

# ---------------------------------------------------------------------------
# Fields sheet: D24 (outputs / output_name / "Possible values") gets the same
# rich text as D19 -- copy/paste so the shared string is reused instead of
# duplicated.
# ---------------------------------------------------------------------------
$cell19.Copy()
$wsFields.Range("D24").PasteSpecial()
$excel.CutCopyMode = $false

# Row 24 grows taller to fit the new text
$wsFields.Rows.Item(24).RowHeight = 151

# ---------------------------------------------------------------------------
# View-state tweaks captured in the diff
# ---------------------------------------------------------------------------
$wsTables.Activate()
$wsTables.Range("C5").Select()

$wsFields.Activate()
$wsFields.Range("C19").Select()
